$d = $word.ActiveDocument

# Locate the target text ("North Central University") robustly via Find.
$target = $d.Content.Duplicate
$found = $target.Find.Execute("North Central University")
if (-not $found) {
    throw "Could not find 'North Central University'"
}

$start = $target.Start

# Offsets within the found text "North Central University":
#   North(0-5) space(5-6) C(6-7) entral University(7-24)
$northEnd = $start + 5   # position right after "North"
$spaceEnd = $start + 6   # position right after the space following "North"
$cEnd     = $start + 7   # position right after the capital "C"

# 1) Lower-case the "C" -> "c" (school's name is one word: Northcentral)
$cRange = $d.Range($northEnd + 1, $cEnd)
$cRange.Text = "c"

# 2) Remove the space between "North" and "c"
$spaceRange = $d.Range($northEnd, $spaceEnd)
$spaceRange.Delete()

# After the two edits above the paragraph reads "Northcentral University"
# and the (now lower-case) "c" sits at [$northEnd, $northEnd + 1).

# 3) Split the single run into three runs: "North" | "c" | "entral University"
#    Dropping a bookmark at each split point forces a run boundary there;
#    deleting the bookmark again removes the bookmark markup but leaves the
#    run split in place.
$split1 = $d.Range($northEnd, $northEnd)
$d.Bookmarks.Add("ntc_split_1", $split1)

$split2 = $d.Range($northEnd + 1, $northEnd + 1)
$d.Bookmarks.Add("ntc_split_2", $split2)

$d.Bookmarks.Item("ntc_split_1").Delete()
$d.Bookmarks.Item("ntc_split_2").Delete()
